$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "68.032.77"
$ws.Cells.Item(2, 5).Value = "  -0.66%  "

$ws.Cells.Item(3, 4).Value = "2.434.96"
$ws.Cells.Item(3, 5).Value = "  -0.36%  "

$ws.Cells.Item(4, 5).Value = "  +0.03%  "

$ws.Cells.Item(5, 4).Value = "'553.98"
$ws.Cells.Item(5, 5).Value = "  -0.22%  "

$ws.Cells.Item(6, 4).Value = "'160.12"
$ws.Cells.Item(6, 5).Value = "  -0.41%  "

$ws.Cells.Item(7, 5).Value = "  +0.08%  "

$ws.Cells.Item(8, 5).Value = "  +2.16%  "

$ws.Cells.Item(9, 4).Value = "'0.159"
$ws.Cells.Item(9, 5).Value = "  +7.75%  "

$ws.Cells.Item(10, 5).Value = "  -0.41%  "

$ws.Cells.Item(11, 5).Value = "  -1.10%  "

$ws.Cells.Item(12, 5).Value = "  +0.42%  "

$ws.Cells.Item(13, 4).Value = "68.501.66"
$ws.Cells.Item(13, 5).Value = "  +0.20%  "

$ws.Cells.Item(14, 5).Value = "  +1.34%  "

$ws.Cells.Item(15, 4).Value = "'23.03"
$ws.Cells.Item(15, 5).Value = "  -0.86%  "

$ws.Cells.Item(16, 5).Value = "  -2.99%  "

$ws.Cells.Item(17, 4).Value = "'333.77"
$ws.Cells.Item(17, 5).Value = "  -1.43%  "

$ws.Cells.Item(18, 5).Value = "  -1.96%  "

$ws.Cells.Item(19, 5).Value = "  +0.13%  "

$ws.Cells.Item(20, 5).Value = "  +0.07%  "

$ws.Cells.Item(21, 5).Value = "  +0.30%  "

$ws.Cells.Item(22, 4).Value = "'66.19"
$ws.Cells.Item(22, 5).Value = "  +0.03%  "

$ws.Cells.Item(23, 5).Value = "  +0.31%  "

$ws.Cells.Item(24, 4).Value = "'8.11"
$ws.Cells.Item(24, 5).Value = "  +0.76%  "

$ws.Cells.Item(25, 5).Value = "  +0.05%  "

$ws.Cells.Item(26, 5).Value = "  +0.04%  "

$ws.Cells.Item(27, 4).Value = "'0.999"
$ws.Cells.Item(27, 5).Value = "  +0.00%  "

$ws.Cells.Item(28, 4).Value = "'419.07"
$ws.Cells.Item(28, 5).Value = "  -3.02%  "

$ws.Cells.Item(29, 5).Value = "  +1.20%  "

$ws.Cells.Item(30, 5).Value = "  +0.02%  "

$ws.Cells.Item(31, 4).Value = "'160.60"
$ws.Cells.Item(31, 5).Value = "  +2.94%  "

$ws.Cells.Item(32, 5).Value = "  -0.40%  "

$ws.Cells.Item(33, 5).Value = "  +0.02%  "

$ws.Cells.Item(34, 4).Value = "'17.86"
$ws.Cells.Item(34, 5).Value = "  +0.63%  "

$ws.Cells.Item(35, 5).Value = "  -3.36%  "

$ws.Cells.Item(36, 5).Value = "  -2.14%  "

$ws.Cells.Item(37, 4).Value = "'4.27"
$ws.Cells.Item(37, 5).Value = "  -2.67%  "

$ws.Cells.Item(38, 4).Value = "'1.47"
$ws.Cells.Item(38, 5).Value = "  +1.27%  "

$ws.Cells.Item(39, 5).Value = "  -0.78%  "

$ws.Cells.Item(40, 5).Value = "  -1.26%  "

$ws.Cells.Item(41, 4).Value = "'3.33"
$ws.Cells.Item(41, 5).Value = "  +0.48%  "

$ws.Cells.Item(42, 4).Value = "'129.59"
$ws.Cells.Item(42, 5).Value = "  -1.39%  "

$ws.Cells.Item(43, 4).Value = "'0.0708"
$ws.Cells.Item(43, 5).Value = "  -0.34%  "

$ws.Cells.Item(44, 4).Value = "'0.479"
$ws.Cells.Item(44, 5).Value = "  -0.05%  "

$ws.Cells.Item(45, 5).Value = "  -0.19%  "

$ws.Cells.Item(46, 4).Value = "'0.0913"
$ws.Cells.Item(46, 5).Value = "  +1.10%  "

$ws.Cells.Item(47, 5).Value = "  +0.53%  "

$ws.Cells.Item(48, 5).Value = "  -5.82%  "

$ws.Cells.Item(49, 4).Value = "'16.56"
$ws.Cells.Item(49, 5).Value = "  -0.97%  "

$ws.Cells.Item(50, 4).Value = "0.0{0}0204" -f [char]0x2086
$ws.Cells.Item(50, 5).Value = "  +4.43%  "

$ws.Cells.Item(51, 5).Value = "  +0.67%  "
